$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.960.02"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.844.80"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'232.35"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'40.77"
$ws.Range("E8").Value = "  +4.33%  "
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").Value = "'0.0984"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "2.113.38"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'11.43"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").Value = "1.842.00"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "34.991.69"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'69.94"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").Value = "'240.26"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "'172.51"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'7.82"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("D29").Value = "'1.64"
$ws.Range("E29").Value = "  +6.49%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'3.96"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +23.66%  "
$ws.Range("E35").Value = "  +11.20%  "
$ws.Range("D36").Value = "'0.751"
$ws.Range("E36").Value = "  +10.48%  "
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("E38").Value = "  +12.01%  "
$ws.Range("D39").Value = "'89.93"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "1.346.58"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").Value = "'14.65"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").Value = "'0.0530"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").Value = "2.030.87"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "'3.43"
$ws.Range("E49").Value = "  +17.73%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -0.39%  "
